$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet has three side-by-side "Code Freq / O-Scope Freq" blocks
# (A:B, D:E, G:H). A new "Jittery?" column is being added after the
# O-Scope Freq column of each block, so insert one blank column in
# front of the 2nd block (old D) and one in front of the 3rd block
# (which, after the first insert, is the new H) - this pushes the
# existing data from D/E -> E/F and G/H -> I/J.

# Remember the widths that were already set on the O-Scope Freq columns
# so they can be restored/reapplied after the insert shifts things.
$widthB = $ws.Columns("B:B").ColumnWidth
$widthE = $ws.Columns("E:E").ColumnWidth
$widthH = $ws.Columns("H:H").ColumnWidth

$ws.Columns("D:D").Insert()
$ws.Columns("H:H").Insert()

# Label the three new "Jittery?" columns.
$ws.Range("C1").Value = "Jittery?"
$ws.Range("G1").Value = "Jittery?"
$ws.Range("K1").Value = "Jittery?"

# E and H used to carry custom widths (as the O-Scope Freq columns of
# the 2nd/3rd blocks); re-apply those same widths now that E and H hold
# different content, and give the new O-Scope Freq columns (F, J) the
# width that column B (the 1st block's O-Scope Freq column) already has.
$ws.Columns("E:E").ColumnWidth = $widthE
$ws.Columns("F:F").ColumnWidth = $widthB
$ws.Columns("H:H").ColumnWidth = $widthH
$ws.Columns("J:J").ColumnWidth = $widthE

# The user's last selection before saving was J10.
$ws.Range("J10").Select()
